$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.416.17"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "3.458.01"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'588.94"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'176.72"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").Value = "'0.614"
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "3.454.54"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").Value = "'0.134"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").Value = "4.057.00"
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "'30.27"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "66.295.93"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "'0.0000173"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "3.458.29"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("D20").Value = "'13.82"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("D21").Value = "'373.97"
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("D22").Value = "'7.62"
$ws.Range("E22").Value = "  -3.42%  "
$ws.Range("D23").Value = "'73.28"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("D27").Value = "'9.91"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'5.88"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").Value = "'23.73"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'7.04"
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("E35").Value = "  -7.00%  "
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("D37").Value = "'160.56"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "'0.885"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "'28.46"
$ws.Range("E39").Value = "  -4.30%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").Value = "2.766.75"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").Value = "'6.46"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("D46").Value = "'25.27"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").Value = "'339.63"
$ws.Range("E47").Value = "  +3.39%  "
$ws.Range("D48").Value = "'40.05"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("D49").Value = "'0.0293"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'31.62"
$ws.Range("E51").Value = "  -0.62%  "
